$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 233, shifting existing rows 233:296 down to 234:297.
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with the new data record.
$ws.Cells.Item(233, 1).Value2 = 4
$ws.Cells.Item(233, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(233, 3).Value2 = "Los Lagos"
$ws.Cells.Item(233, 4).Value2 = 44943
$ws.Cells.Item(233, 5).Value2 = 10
$ws.Cells.Item(233, 6).Value2 = 100112039
$ws.Cells.Item(233, 7).Value2 = "Ciboulette"
$ws.Cells.Item(233, 8).Value2 = "Sin especificar"
$ws.Cells.Item(233, 9).Value2 = "Primera"
$ws.Cells.Item(233, 10).Value2 = 240
$ws.Cells.Item(233, 11).Value2 = 3000
$ws.Cells.Item(233, 12).Value2 = 3500
$ws.Cells.Item(233, 13).Value2 = 3250
$ws.Cells.Item(233, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(233, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(233, 16).Value2 = 1083
$ws.Cells.Item(233, 17).Value2 = 3
$ws.Cells.Item(233, 18).Value2 = "Hortaliza"
